$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 200, shifting the existing
# rows 200-233 down to 201-234 (weekly refresh adding one new price record).
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with the new week's record.
$ws.Range("A200").Value = 11
$ws.Range("B200").Value = "Vega Monumental Concepción"
$ws.Range("C200").Value = "Bíobío"
$ws.Range("D200").Value = 45127
$ws.Range("E200").Value = 8
$ws.Range("F200").Value = "Fruta"
$ws.Range("G200").Value = 100102
$ws.Range("H200").Value = "Cítricos"
$ws.Range("I200").Value = 100102004
$ws.Range("J200").Value = "Mandarina"
$ws.Range("K200").Value = "Clemenuless"
$ws.Range("L200").Value = "Primera"
$ws.Range("M200").Value = 200
$ws.Range("N200").Value = 9000
$ws.Range("O200").Value = 9500
$ws.Range("P200").Value = 9250
$ws.Range("Q200").Value = "$/bandeja 10 kilos"
$ws.Range("R200").Value = "Provincia de Limarí"
$ws.Range("S200").Value = 925
$ws.Range("T200").Value = 10
